$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the duplicate X:AG columns (old sheet had 32 HKL-label columns,
#    new one has 21; X:AG mirrored N:W and is no longer needed).
$ws.Range("X1:AG19").EntireColumn.Delete()

# 2) Reorder the HKL header labels in row 2 (C2:K2). L2/M2 and N2:W2 keep
#    their existing labels.
$ws.Range("C2").Value = "[4, 2, 0]"
$ws.Range("D2").Value = "[4, 0, 0]"
$ws.Range("E2").Value = "[2, 0, 0]"
$ws.Range("F2").Value = "[2, 2, 0]"
$ws.Range("G2").Value = "[3, 3, 3]"
$ws.Range("H2").Value = "[1, 1, 1]"
$ws.Range("I2").Value = "[2, 2, 2]"
$ws.Range("J2").Value = "[3, 3, 1]"
$ws.Range("K2").Value = "[3, 1, 1]"

# 3) Reorder the HKL data columns C:K for the rows that keep their scheme
#    (rows 3-15: BT8Hex/Spiral/Offset); also re-stamp L:W with the values
#    from the regenerated simulation (mostly identical, last-digit noise).
$ws.Range("C3").Value = 1.008778341630243
$ws.Range("D3").Value = 1.038098124455384
$ws.Range("E3").Value = 1.038098124455384
$ws.Range("F3").Value = 0.9991323977173241
$ws.Range("G3").Value = 0.9718890071522363
$ws.Range("H3").Value = 0.9718890071522363
$ws.Range("I3").Value = 0.9718890071522363
$ws.Range("J3").Value = 0.9896895619473145
$ws.Range("K3").Value = 1.004262683892209
$ws.Range("L3").Value = 0.9869348789474791
$ws.Range("M3").Value = 1.022213094466405
$ws.Range("N3").Value = 0.9718890071522363
$ws.Range("O3").Value = 0.9991323977173241
$ws.Range("P3").Value = 1.018615261086354
$ws.Range("Q3").Value = 1.001697540804766
$ws.Range("R3").Value = 1.003039843108315
$ws.Range("S3").Value = 1.013831068688305
$ws.Range("T3").Value = 1.003039843108315
$ws.Range("U3").Value = 1.003345553304288
$ws.Range("V3").Value = 0.9970542440738777
$ws.Range("W3").Value = 1.002624761276074

$ws.Range("C4").Value = 1.017024146827711
$ws.Range("D4").Value = 1.072948661355545
$ws.Range("E4").Value = 1.072948661355545
$ws.Range("F4").Value = 0.997907526360858
$ws.Range("G4").Value = 0.9455123443089511
$ws.Range("H4").Value = 0.9455123443089511
$ws.Range("I4").Value = 0.9455123443089511
$ws.Range("J4").Value = 0.9798891519422757
$ws.Range("K4").Value = 1.008630370789052
$ws.Range("L4").Value = 0.975086319551327
$ws.Range("M4").Value = 1.042902422397654
$ws.Range("N4").Value = 0.9455123443089511
$ws.Range("O4").Value = 0.997907526360858
$ws.Range("P4").Value = 1.035428093858201
$ws.Range("Q4").Value = 1.003268948574955
$ws.Range("R4").Value = 1.005456177341785
$ws.Range("S4").Value = 1.026495519501818
$ws.Range("T4").Value = 1.005456177341785
$ws.Range("U4").Value = 1.006249725703601
$ws.Range("V4").Value = 0.9941022494246713
$ws.Range("W4").Value = 1.004987617941671

$ws.Range("C5").Value = 1.032059379990903
$ws.Range("D5").Value = 1.141477102769183
$ws.Range("E5").Value = 1.141477102769183
$ws.Range("F5").Value = 0.9957370190239905
$ws.Range("G5").Value = 0.8969173518638881
$ws.Range("H5").Value = 0.8969173518638881
$ws.Range("I5").Value = 0.8969173518638881
$ws.Range("J5").Value = 0.9615188317563654
$ws.Range("K5").Value = 1.015946930822984
$ws.Range("L5").Value = 0.95208297331382
$ws.Range("M5").Value = 1.082407587582736
$ws.Range("N5").Value = 0.8969173518638881
$ws.Range("O5").Value = 0.9957370190239905
$ws.Range("P5").Value = 1.068607060896587
$ws.Range("Q5").Value = 1.005841974923487
$ws.Range("R5").Value = 1.011377157885687
$ws.Range("S5").Value = 1.051053684205386
$ws.Range("T5").Value = 1.011377157885687
$ws.Range("U5").Value = 1.012519601120011
$ws.Range("V5").Value = 0.9893991512687869
$ws.Range("W5").Value = 1.009768397140484

$ws.Range("C6").Value = 1.047397828025148
$ws.Range("D6").Value = 1.208184877687989
$ws.Range("E6").Value = 1.208184877687989
$ws.Range("F6").Value = 0.993467589857899
$ws.Range("G6").Value = 0.8503033370542593
$ws.Range("H6").Value = 0.8503033370542593
$ws.Range("I6").Value = 0.8503033370542593
$ws.Range("J6").Value = 0.9430601767585016
$ws.Range("K6").Value = 1.022855564121425
$ws.Range("L6").Value = 0.9294372372200028
$ws.Range("M6").Value = 1.120761916832289
$ws.Range("N6").Value = 0.8503033370542593
$ws.Range("O6").Value = 0.993467589857899
$ws.Range("P6").Value = 1.100826233772944
$ws.Range("Q6").Value = 1.008161576989662
$ws.Range("R6").Value = 1.017318601533383
$ws.Range("S6").Value = 1.074836010555771
$ws.Range("T6").Value = 1.017318601533382
$ws.Range("U6").Value = 1.018702842180393
$ws.Range("V6").Value = 0.9850229411551663
$ws.Range("W6").Value = 1.014433565944689

$ws.Range("C7").Value = 1.000218758635835
$ws.Range("D7").Value = 1.003618233355279
$ws.Range("E7").Value = 1.003618233355279
$ws.Range("F7").Value = 0.9990234812033179
$ws.Range("G7").Value = 0.9987831212398735
$ws.Range("H7").Value = 0.9987831212398735
$ws.Range("I7").Value = 0.9987831212398735
$ws.Range("J7").Value = 0.9989536758651959
$ws.Range("K7").Value = 1.000398622548761
$ws.Range("L7").Value = 0.9993153381581527
$ws.Range("M7").Value = 1.001947133459606
$ws.Range("N7").Value = 0.9987831212398735
$ws.Range("O7").Value = 0.9990234812033179
$ws.Range("P7").Value = 1.001320857279298
$ws.Range("Q7").Value = 0.9997110518760396
$ws.Range("R7").Value = 1.000474945266157
$ws.Range("S7").Value = 1.001013445702453
$ws.Range("T7").Value = 1.000474945266157
$ws.Range("U7").Value = 1.000455864586808
$ws.Range("V7").Value = 1.000121315917421
$ws.Range("W7").Value = 1.000282295558253

$ws.Range("C8").Value = 1.000653040086157
$ws.Range("D8").Value = 1.009776360203636
$ws.Range("E8").Value = 1.009776360203636
$ws.Range("F8").Value = 0.9973512063936659
$ws.Range("G8").Value = 0.9966175327164225
$ws.Range("H8").Value = 0.9966175327164225
$ws.Range("I8").Value = 0.9966175327164225
$ws.Range("J8").Value = 0.9971252964540895
$ws.Range("K8").Value = 1.00111847975456
$ws.Range("L8").Value = 0.9981176711792353
$ws.Range("M8").Value = 1.005317756988708
$ws.Range("N8").Value = 0.9966175327164225
$ws.Range("O8").Value = 0.9973512063936659
$ws.Range("P8").Value = 1.003563783298651
$ws.Range("Q8").Value = 0.9992348430741131
$ws.Range("R8").Value = 1.001248366437908
$ws.Range("S8").Value = 1.002748682117287
$ws.Range("T8").Value = 1.001248366437908
$ws.Range("U8").Value = 1.001215894767071
$ws.Range("V8").Value = 1.000296222356941
$ws.Range("W8").Value = 1.000759667972059

$ws.Range("C9").Value = 1.00095334700493
$ws.Range("D9").Value = 1.014093610266912
$ws.Range("E9").Value = 1.014093610266912
$ws.Range("F9").Value = 0.996120008700757
$ws.Range("G9").Value = 0.9950297917149249
$ws.Range("H9").Value = 0.9950297917149249
$ws.Range("I9").Value = 0.9950297917149249
$ws.Range("J9").Value = 0.9958186028064048
$ws.Range("K9").Value = 1.001675560609047
$ws.Range("L9").Value = 0.9973155417667406
$ws.Range("M9").Value = 1.007705340839672
$ws.Range("N9").Value = 0.9950297917149249
$ws.Range("O9").Value = 0.996120008700757
$ws.Range("P9").Value = 1.005106809483834
$ws.Range("Q9").Value = 0.9988977846549021
$ws.Range("R9").Value = 1.001747803560865
$ws.Range("S9").Value = 1.003963059858905
$ws.Range("T9").Value = 1.001747803560865
$ws.Range("U9").Value = 1.00172974282291
$ws.Range("V9").Value = 1.000389752601313
$ws.Range("W9").Value = 1.001088975463674

$ws.Range("C10").Value = 1.002015240045356
$ws.Range("D10").Value = 1.02997296244758
$ws.Range("E10").Value = 1.02997296244758
$ws.Range("F10").Value = 0.9918647936779412
$ws.Range("G10").Value = 0.98936624571685
$ws.Range("H10").Value = 0.98936624571685
$ws.Range("I10").Value = 0.98936624571685
$ws.Range("J10").Value = 0.9911704266759731
$ws.Range("K10").Value = 1.003525513258995
$ws.Range("L10").Value = 0.9942770351226389
$ws.Range("M10").Value = 1.0163389991625
$ws.Range("N10").Value = 0.9893662457168499
$ws.Range("O10").Value = 0.9918647936779412
$ws.Range("P10").Value = 1.01091887806276
$ws.Range("Q10").Value = 0.997695153468468
$ws.Range("R10").Value = 1.00373466728079
$ws.Range("S10").Value = 1.008454423128172
$ws.Range("T10").Value = 1.00373466728079
$ws.Range("U10").Value = 1.003682378775341
$ws.Range("V10").Value = 1.000819152163643
$ws.Range("W10").Value = 1.002316402013479

$ws.Range("C11").Value = 1.003370062044892
$ws.Range("D11").Value = 1.051256125319738
$ws.Range("E11").Value = 1.051256125319738
$ws.Range("F11").Value = 0.98597394128689
$ws.Range("G11").Value = 0.9821935340076484
$ws.Range("H11").Value = 0.9821935340076484
$ws.Range("I11").Value = 0.9821935340076484
$ws.Range("J11").Value = 0.9848885506130653
$ws.Range("K11").Value = 1.005967386993258
$ws.Range("L11").Value = 0.9902387075723729
$ws.Range("M11").Value = 1.02791422250914
$ws.Range("N11").Value = 0.9821935340076484
$ws.Range("O11").Value = 0.98597394128689
$ws.Range("P11").Value = 1.018615033303314
$ws.Range("Q11").Value = 0.9959706641400738
$ws.Range("R11").Value = 1.006474533538092
$ws.Range("S11").Value = 1.014399151199962
$ws.Range("T11").Value = 1.006474533538092
$ws.Range("U11").Value = 1.006347746901884
$ws.Range("V11").Value = 1.001516904323037
$ws.Range("W11").Value = 1.003975316293376

$ws.Range("C12").Value = 0.8745368999227625
$ws.Range("D12").Value = 0.6132311431620122
$ws.Range("E12").Value = 0.6132311431620122
$ws.Range("F12").Value = 0.9489468280812738
$ws.Range("G12").Value = 1.391676314260136
$ws.Range("H12").Value = 1.391676314260136
$ws.Range("I12").Value = 1.391676314260136
$ws.Range("J12").Value = 1.097446714224897
$ws.Range("K12").Value = 0.9525803801333577
$ws.Range("L12").Value = 1.16504839110238
$ws.Range("M12").Value = 0.7638635791204468
$ws.Range("N12").Value = 1.391676314260136
$ws.Range("O12").Value = 0.9489468280812738
$ws.Range("P12").Value = 0.781088985621643
$ws.Range("Q12").Value = 0.9507636041073158
$ws.Range("R12").Value = 0.9846180951678072
$ws.Range("S12").Value = 0.8382527837922146
$ws.Range("T12").Value = 0.9846180951678072
$ws.Range("U12").Value = 0.9766086664091949
$ws.Range("V12").Value = 1.059622195979383
$ws.Range("W12").Value = 0.9759162812509081

$ws.Range("C13").Value = 0.9921446998649259
$ws.Range("D13").Value = 0.9208543484117759
$ws.Range("E13").Value = 0.9208543484117759
$ws.Range("F13").Value = 1.114500948424145
$ws.Range("G13").Value = 1.044247981991259
$ws.Range("H13").Value = 1.044247981991259
$ws.Range("I13").Value = 1.044247981991259
$ws.Range("J13").Value = 1.074483415934369
$ws.Range("K13").Value = 0.9330467897408178
$ws.Range("L13").Value = 0.9797230111250645
$ws.Range("M13").Value = 0.9170928555211185
$ws.Range("N13").Value = 1.044247981991259
$ws.Range("O13").Value = 1.114500948424145
$ws.Range("P13").Value = 1.017677648417961
$ws.Range("Q13").Value = 1.023773869082481
$ws.Range("R13").Value = 1.026534426275727
$ws.Range("S13").Value = 0.9894673621922463
$ws.Range("T13").Value = 1.026534426275727
$ws.Range("U13").Value = 1.003162517142
$ws.Range("V13").Value = 1.011379610111852
$ws.Range("W13").Value = 0.9970117563766845

$ws.Range("C14").Value = 1.097375869095959
$ws.Range("D14").Value = 0.6474123980439753
$ws.Range("E14").Value = 0.6474123980439753
$ws.Range("F14").Value = 1.383247727733637
$ws.Range("G14").Value = 0.7695134981280533
$ws.Range("H14").Value = 0.7695134981280533
$ws.Range("I14").Value = 0.7695134981280533
$ws.Range("J14").Value = 1.176299178148569
$ws.Range("K14").Value = 0.9210545075827427
$ws.Range("L14").Value = 0.9294315679662116
$ws.Range("M14").Value = 0.8036241904562451
$ws.Range("N14").Value = 0.7695134981280533
$ws.Range("O14").Value = 1.383247727733637
$ws.Range("P14").Value = 1.015330062888806
$ws.Range("Q14").Value = 1.15215111765819
$ws.Range("R14").Value = 0.9333912079685552
$ws.Range("S14").Value = 0.983904877786785
$ws.Range("T14").Value = 0.9333912079685552
$ws.Range("U14").Value = 0.9303070328721021
$ws.Range("V14").Value = 0.8981483259232924
$ws.Range("W14").Value = 0.9659948671444241

$ws.Range("C15").Value = 0.9361305276650942
$ws.Range("D15").Value = 1.013927747195857
$ws.Range("E15").Value = 1.013927747195857
$ws.Range("F15").Value = 0.8574694106388464
$ws.Range("G15").Value = 1.165039530280905
$ws.Range("H15").Value = 1.165039530280905
$ws.Range("I15").Value = 1.165039530280905
$ws.Range("J15").Value = 0.9650377884299971
$ws.Range("K15").Value = 1.020185115454781
$ws.Range("L15").Value = 1.069688241736297
$ws.Range("M15").Value = 1.00572304026691
$ws.Range("N15").Value = 1.165039530280905
$ws.Range("O15").Value = 0.8574694106388464
$ws.Range("P15").Value = 0.9356985789173518
$ws.Range("Q15").Value = 0.9388272630468135
$ws.Range("R15").Value = 1.012145562705203
$ws.Range("S15").Value = 0.9638607577631614
$ws.Range("T15").Value = 1.012145562705203
$ws.Range("U15").Value = 1.014155450892597
$ws.Range("V15").Value = 1.044332266770259
$ws.Range("W15").Value = 1.004150175208586

# 4) Replace rows 16-19 with the new Holden scheme (label + full data row).
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("C16").Value = 1.167276980881863
$ws.Range("D16").Value = 1.738533591951291
$ws.Range("E16").Value = 1.738533591951291
$ws.Range("F16").Value = 0.9629589977788179
$ws.Range("G16").Value = 0.459752765546561
$ws.Range("H16").Value = 0.459752765546561
$ws.Range("I16").Value = 0.459752765546561
$ws.Range("J16").Value = 0.791282979177104
$ws.Range("K16").Value = 1.092449638973112
$ws.Range("L16").Value = 0.755221057606827
$ws.Range("M16").Value = 1.437002068385266
$ws.Range("N16").Value = 0.459752765546561
$ws.Range("O16").Value = 0.9629589977788179
$ws.Range("P16").Value = 1.350746294865054
$ws.Range("Q16").Value = 1.027704318375965
$ws.Range("R16").Value = 1.05374845175889
$ws.Range("S16").Value = 1.26464740956774
$ws.Range("T16").Value = 1.05374845175889
$ws.Range("U16").Value = 1.063423748562446
$ws.Range("V16").Value = 0.9426895519592687
$ws.Range("W16").Value = 1.050559760037605

$ws.Range("B17").Value = "Holden5"
$ws.Range("C17").Value = 1.124335497661176
$ws.Range("D17").Value = 1.620490279891411
$ws.Range("E17").Value = 1.620490279891411
$ws.Range("F17").Value = 0.9394744019637108
$ws.Range("G17").Value = 0.5857047570070888
$ws.Range("H17").Value = 0.5857047570070888
$ws.Range("I17").Value = 0.5857047570070888
$ws.Range("J17").Value = 0.8197787156013153
$ws.Range("K17").Value = 1.080348720686533
$ws.Range("L17").Value = 0.8108199619374276
$ws.Range("M17").Value = 1.365148103960654
$ws.Range("N17").Value = 0.5857047570070888
$ws.Range("O17").Value = 0.9394744019637108
$ws.Range("P17").Value = 1.279982340927561
$ws.Range("Q17").Value = 1.009911561325122
$ws.Range("R17").Value = 1.048556479620737
$ws.Range("S17").Value = 1.213437800847218
$ws.Range("T17").Value = 1.048556479620737
$ws.Range("U17").Value = 1.056504539887186
$ws.Range("V17").Value = 0.9623445833111666
$ws.Range("W17").Value = 1.043262554838664

$ws.Range("B18").Value = "Holden10"
$ws.Range("C18").Value = 1.038657238484366
$ws.Range("D18").Value = 1.383256341173908
$ws.Range("E18").Value = 1.383256341173908
$ws.Range("F18").Value = 0.893189163622372
$ws.Range("G18").Value = 0.8373413992192996
$ws.Range("H18").Value = 0.8373413992192996
$ws.Range("I18").Value = 0.8373413992192996
$ws.Range("J18").Value = 0.8771192296639893
$ws.Range("K18").Value = 1.056038767607439
$ws.Range("L18").Value = 0.9219960476057456
$ws.Range("M18").Value = 1.220906361109976
$ws.Range("N18").Value = 0.8373413992192996
$ws.Range("O18").Value = 0.893189163622372
$ws.Range("P18").Value = 1.13822275239814
$ws.Range("Q18").Value = 0.9746139656149055
$ws.Range("R18").Value = 1.037928968005193
$ws.Range("S18").Value = 1.11082809080124
$ws.Range("T18").Value = 1.037928968005193
$ws.Range("U18").Value = 1.042456417905755
$ws.Range("V18").Value = 1.001433414168464
$ws.Range("W18").Value = 1.028563068560887

$ws.Range("B19").Value = "Holden15"
$ws.Range("C19").Value = 1.032744548140657
$ws.Range("D19").Value = 1.448363226255342
$ws.Range("E19").Value = 1.448363226255342
$ws.Range("F19").Value = 0.8849019364999398
$ws.Range("G19").Value = 0.8404673446153117
$ws.Range("H19").Value = 0.8404673446153117
$ws.Range("I19").Value = 0.8404673446153117
$ws.Range("J19").Value = 0.8698033363328729
$ws.Range("K19").Value = 1.049285089024483
$ws.Range("L19").Value = 0.9094462423586924
$ws.Range("M19").Value = 1.243947219971435
$ws.Range("N19").Value = 0.8404673446153117
$ws.Range("O19").Value = 0.8849019364999398
$ws.Range("P19").Value = 1.166632581377641
$ws.Range("Q19").Value = 0.9670935127622116
$ws.Range("R19").Value = 1.057910835790198
$ws.Range("S19").Value = 1.127516750593255
$ws.Range("T19").Value = 1.057910835790198
$ws.Range("U19").Value = 1.055754399098769
$ws.Range("V19").Value = 1.012696988202078
$ws.Range("W19").Value = 1.034869867899842

# 5) Append rows 20-23 with the HexGrid data that used to live in rows 16-19
#    (same permutation applied to C:K; L:W carried over unchanged).
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "HexGrid-90degTilt2.5degRes"
$ws.Range("C20").Value = 0.9997845127872559
$ws.Range("D20").Value = 1.000579898252708
$ws.Range("E20").Value = 1.000579898252708
$ws.Range("F20").Value = 0.9995246807456281
$ws.Range("G20").Value = 1.000399194364444
$ws.Range("H20").Value = 1.000399194364444
$ws.Range("I20").Value = 1.000399194364444
$ws.Range("J20").Value = 0.999833183913002
$ws.Range("K20").Value = 1.0000279668196
$ws.Range("L20").Value = 1.000108099830392
$ws.Range("M20").Value = 1.000211242463924
$ws.Range("N20").Value = 1.000399194364444
$ws.Range("O20").Value = 0.9995246807456281
$ws.Range("P20").Value = 1.000052289499168
$ws.Range("Q20").Value = 0.9997763237826138
$ws.Range("R20").Value = 1.00016792445426
$ws.Range("S20").Value = 1.000044181939312
$ws.Range("T20").Value = 1.00016792445426
$ws.Range("U20").Value = 1.000132935045595
$ws.Range("V20").Value = 1.000186186909365
$ws.Range("W20").Value = 1.000058597397119

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C21").Value = 1.000438000469342
$ws.Range("D21").Value = 1.003839306369855
$ws.Range("E21").Value = 1.003839306369855
$ws.Range("F21").Value = 0.9994168266002679
$ws.Range("G21").Value = 0.998128435043112
$ws.Range("H21").Value = 0.998128435043112
$ws.Range("I21").Value = 0.998128435043112
$ws.Range("J21").Value = 0.998989343418796
$ws.Range("K21").Value = 1.00038119888568
$ws.Range("L21").Value = 0.9990310435735678
$ws.Range("M21").Value = 1.002098984670025
$ws.Range("N21").Value = 0.998128435043112
$ws.Range("O21").Value = 0.9994168266002679
$ws.Range("P21").Value = 1.001628066485061
$ws.Range("Q21").Value = 0.9998990127429738
$ws.Range("R21").Value = 1.000461522671078
$ws.Range("S21").Value = 1.001212443951934
$ws.Range("T21").Value = 1.000461522671078
$ws.Range("U21").Value = 1.000441441724728
$ws.Range("V21").Value = 0.9999788403884052
$ws.Range("W21").Value = 1.00029039237883

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "HexGrid-90degTilt10degRes"
$ws.Range("C22").Value = 1.002545680883234
$ws.Range("D22").Value = 1.011839892986081
$ws.Range("E22").Value = 1.011839892986081
$ws.Range("F22").Value = 0.9988825685582906
$ws.Range("G22").Value = 0.9917039195673589
$ws.Range("H22").Value = 0.9917039195673589
$ws.Range("I22").Value = 0.9917039195673589
$ws.Range("J22").Value = 0.9965067294468001
$ws.Range("K22").Value = 1.001574625979797
$ws.Range("L22").Value = 0.9964016031586036
$ws.Range("M22").Value = 1.006919911905917
$ws.Range("N22").Value = 0.9917039195673589
$ws.Range("O22").Value = 0.9988825685582906
$ws.Range("P22").Value = 1.005361230772186
$ws.Range("Q22").Value = 1.000228597269044
$ws.Range("R22").Value = 1.00080879370391
$ws.Range("S22").Value = 1.004099029174723
$ws.Range("T22").Value = 1.00080879370391
$ws.Range("U22").Value = 1.001000251772882
$ws.Range("V22").Value = 0.9991409853317773
$ws.Range("W22").Value = 1.00079686656076

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "HexGrid-90degTilt15degRes"
$ws.Range("C23").Value = 1.00651111731119
$ws.Range("D23").Value = 1.027017926107413
$ws.Range("E23").Value = 1.027017926107413
$ws.Range("F23").Value = 0.9993529746897861
$ws.Range("G23").Value = 0.9792462999134202
$ws.Range("H23").Value = 0.9792462999134202
$ws.Range("I23").Value = 0.9792462999134202
$ws.Range("J23").Value = 0.992505511108699
$ws.Range("K23").Value = 1.003283701559858
$ws.Range("L23").Value = 0.9907141891232444
$ws.Range("M23").Value = 1.015927470897878
$ws.Range("N23").Value = 0.9792462999134202
$ws.Range("O23").Value = 0.9993529746897861
$ws.Range("P23").Value = 1.0131854503986
$ws.Range("Q23").Value = 1.001318338124822
$ws.Range("R23").Value = 1.001872400236873
$ws.Range("S23").Value = 1.009884867452352
$ws.Range("T23").Value = 1.001872400236873
$ws.Range("U23").Value = 1.002225225567619
$ws.Range("V23").Value = 0.9976294404367796
$ws.Range("W23").Value = 1.001819898838936

# 6) Column A on the new rows needs the same bordered/bold/centered style as
#    the rest of column A; copy formats only (keeps styles.xml unchanged).
$ws.Range("A3").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

